# add new feeCal function.
#
# Populates the new "PROCESS_DATE" column (I2:I7) with the same dates as
# the existing TRADE_DATE / VALN_DATE columns (E and G), using the same
# centered, custom "yyyy-mm-dd" date format already used by those columns,
# then moves the active selection onto the newly filled range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Same serial date values already present in columns E (TRADE_DATE) and
# G (VALN_DATE) for rows 2-7: 2017-05-01 .. 2017-10-01.
$processDates = @(42856, 42887, 42917, 42948, 42979, 43009)

for ($i = 0; $i -lt $processDates.Length; $i++) {
    $row = 2 + $i
    $cell = $ws.Cells.Item($row, 9)   # column I

    # Apply alignment/format in the same order used elsewhere in the sheet
    # so the cell lands on the existing shared style (center + custom date
    # format) instead of spawning a brand new one.
    $cell.HorizontalAlignment = -4108  # xlCenter
    $cell.NumberFormat = "yyyy\-mm\-dd;@"
    $cell.Value = $processDates[$i]
}

$ws.Range("I2:I7").Select()
